$d = $word.ActiveDocument

# The first paragraph currently reads "UC1 - tjek kreditværdighed".
# Prepend a new, separately-formatted run containing "FFS-" right after
# the _GoBack bookmark (and before the existing bold run), so the
# paragraph becomes two runs: "FFS-" + "UC1 - tjek kreditværdighed".
$firstPara = $d.Paragraphs(1)
$startChar = $firstPara.Range.Characters(1)
$startChar.InsertBefore("FFS-")

# Re-apply the same bold / size formatting as the following run. Toggling
# Bold off then on forces the new text to stay in its own <w:r> instead of
# being coalesced into the adjacent identically-formatted run.
$newRun = $d.Range(0, 4)
$newRun.Font.Bold = $false
$newRun.Font.Bold = $true
$newRun.Font.Size = 12
$newRun.Font.SizeBi = 12
